$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the duplicated bold "Play Dragon Auto Chess Slot Game for
#    Free - Review" paragraph that currently sits near the end of the
#    document (right before the italic meta-description paragraph).
#    We must not touch the Heading1 paragraph with the same text at the
#    very start of the document, so we find the *last* paragraph whose
#    text equals the title and remove it (including its paragraph
#    mark).
# ---------------------------------------------------------------------
$titleText = "Play Dragon Auto Chess Slot Game for Free - Review"
$titleParaToDelete = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText) {
        $titleParaToDelete = $p
    }
}
if ($titleParaToDelete -ne $null) {
    $titleParaToDelete.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Replace the text of the italic paragraph (old meta description)
#    with the new "Prompt: ..." copy, keeping the run's italic
#    formatting and the paragraph's leading empty run untouched. We set
#    Range.Text directly (instead of Find/Replace) so literal straight
#    quotes are preserved verbatim.
# ---------------------------------------------------------------------
$oldBlurb = "Find out what we thought of Dragon Auto Chess, an innovative slot game inspired by chess. Play for free and see if you can win big!"
$newPrompt = 'Prompt: Please create an image in a cartoon style featuring a happy Maya warrior with glasses for the game "Dragon Auto Chess". The background should be a castle, and the Maya warrior should be holding a chess piece or a dragon symbol. Use bright and vibrant colors to capture the fun and exciting nature of the game.'

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $oldBlurb) {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Text = $newPrompt
    }
}

# ---------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph at the top of the document. It contains
#    a leading empty run, a bold "Meta description" run, and a plain
#    run with the rest of the sentence.
# ---------------------------------------------------------------------
$heading1 = $d.Paragraphs(1)
$insertRange = $heading1.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaStart = $metaPara.Range.Start
$boldText = "Meta description"
$restText = ": Find out what we thought of Dragon Auto Chess, an innovative slot game inspired by chess. Play for free and see if you can win big!"

$metaPara.Range.InsertBefore($boldText + $restText)

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

$restRange = $d.Range($metaStart + $boldText.Length, $metaStart + $boldText.Length + $restText.Length)
$restRange.Bold = 0
